$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as scraped on 2023-04-06
# D-column (Price) values are numeric-looking strings; prefix with an apostrophe
# so Excel keeps them as text instead of auto-converting to numbers, matching
# the original inlineStr text storage.

# Row 2
$ws.Range("D2").Value = "'28.077.56"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").Value = "'1.871.28"
$ws.Range("E3").Value = "  -2.17%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").Value = "'313.55"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
$ws.Range("E6").Value = "  +0.21%  "

# Row 7
$ws.Range("D7").Value = "'0.5057"
$ws.Range("E7").Value = "  -0.70%  "

# Row 8
$ws.Range("D8").Value = "'0.3842"
$ws.Range("E8").Value = "  -2.12%  "

# Row 9
$ws.Range("D9").Value = "'0.08658"
$ws.Range("E9").Value = "  -6.95%  "

# Row 10
$ws.Range("D10").Value = "'1.115"
$ws.Range("E10").Value = "  -2.25%  "

# Row 11
$ws.Range("D11").Value = "'41.52"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
$ws.Range("D12").Value = "'6.321"
$ws.Range("E12").Value = "  -1.22%  "

# Row 13
$ws.Range("D13").Value = "'20.65"
$ws.Range("E13").Value = "  -1.24%  "

# Row 14
$ws.Range("D14").Value = "'1.868.84"
$ws.Range("E14").Value = "  -1.82%  "

# Row 15
$ws.Range("E15").Value = "  +0.26%  "

# Row 16
$ws.Range("D16").Value = "'7.173"

# Row 17
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = "  -1.78%  "

# Row 18
$ws.Range("D18").Value = "'90.90"
$ws.Range("E18").Value = "  -1.70%  "

# Row 19
$ws.Range("D19").Value = "'0.06631"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20
$ws.Range("D20").Value = "'18.00"
$ws.Range("E20").Value = "  +0.13%  "

# Row 21
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22
$ws.Range("D22").Value = "'6.089"
$ws.Range("E22").Value = "  -2.32%  "

# Row 23
$ws.Range("D23").Value = "'28.106.32"
$ws.Range("E23").Value = "  -0.48%  "

# Row 24
$ws.Range("D24").Value = "'11.42"
$ws.Range("E24").Value = "  -0.40%  "

# Row 25
$ws.Range("D25").Value = "'2.266"
$ws.Range("E25").Value = "  -2.62%  "

# Row 26
$ws.Range("D26").Value = "'2.557"
$ws.Range("E26").Value = "  -1.10%  "

# Row 27
$ws.Range("D27").Value = "'2.082.52"
$ws.Range("E27").Value = "  -1.90%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'157.25"
$ws.Range("E28").Value = "  -0.42%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.72"
$ws.Range("E29").Value = "  -1.89%  "

# Row 30
$ws.Range("D30").Value = "'126.23"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("E31").Value = "  -2.37%  "

# Row 32
$ws.Range("D32").Value = "'1.059"
$ws.Range("E32").Value = "  -3.45%  "

# Row 33
$ws.Range("D33").Value = "'5.593"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34
$ws.Range("D34").Value = "'3.599"
$ws.Range("E34").Value = "  -0.47%  "

# Row 35
$ws.Range("D35").Value = "'9.638"
$ws.Range("E35").Value = "  -0.77%  "

# Row 36
$ws.Range("D36").Value = "'0.02443"
$ws.Range("E36").Value = "  +0.65%  "

# Row 37
$ws.Range("D37").Value = "'0.06574"
$ws.Range("E37").Value = "  -1.50%  "

# Row 38
$ws.Range("D38").Value = "'0.2172"
$ws.Range("E38").Value = "  -1.18%  "

# Row 39
$ws.Range("D39").Value = "'1.204"
$ws.Range("E39").Value = "  -3.73%  "

# Row 40
$ws.Range("D40").Value = "'1.240"
$ws.Range("E40").Value = "  -4.15%  "

# Row 41
$ws.Range("E41").Value = "  +0.50%  "

# Row 42
$ws.Range("D42").Value = "'0.6360"
$ws.Range("E42").Value = "  -1.50%  "

# Row 43
$ws.Range("D43").Value = "'4.889"
$ws.Range("E43").Value = "  -2.25%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.15"
$ws.Range("E44").Value = "  -1.63%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5986"
$ws.Range("E45").Value = "  -1.01%  "

# Row 46
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.280"
$ws.Range("E46").Value = "  -0.21%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.676"
$ws.Range("E47").Value = "  -1.21%  "

# Row 48
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "'1.229"
$ws.Range("E48").Value = "  +3.53%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.986"
$ws.Range("E49").Value = "  -1.53%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'121.25"
$ws.Range("E50").Value = "  -1.78%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'79.99"
$ws.Range("E51").Value = "  +1.82%  "
